# issue #5: add legislator_id, name, date into dataframe
#
# The upstream scraper now stamps every row of the "股票" (stocks) sheet with
# the three identifying columns that already exist as column headers
# elsewhere in the pipeline: date (filing date), legislator_name and
# legislator_id. This mirrors that by adding three trailing columns (H, I, J)
# to the 股票 worksheet: header labels in row 1, and the constant values for
# every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$filingDate = "2011-11-17"
$legislatorName = "吳育昇"
$legislatorId = 1322

# Find the last used data row (column B holds the stock name on every row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# --- Header row -----------------------------------------------------------
$headerRange = $ws.Range("H1:J1")
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Font.Bold = $true

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows --------------------------------------------------------------
$dateRange = $ws.Range("H2:H" + $lastRow)
$dateRange.NumberFormat = "@"
$dateRange.Value = $filingDate
$dateRange.Style = $ws.Range("C2").Style

$nameRange = $ws.Range("I2:I" + $lastRow)
$nameRange.Value = $legislatorName

$idRange = $ws.Range("J2:J" + $lastRow)
$idRange.Value = $legislatorId
